$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (target) order of rows 2..18 after the edit: each player paired
# with their real position/team, rows reordered so that Stephen Curry is
# first, Mark Williams follows Daniel Gafford, and the two Tyrese's move
# down after Keegan Murray.
$data = @(
    @("Stephen Curry",      "PG,SG", "Golden State Warriors"),
    @("Austin Reaves",      "PG,SG", "Los Angeles Lakers"),
    @("Darius Garland",     "PG",    "Cleveland Cavaliers"),
    @("OG Anunoby",         "SF,PF", "New York Knicks"),
    @("Kevin Durant",       "SF,PF", "Phoenix Suns"),
    @("Karl-Anthony Towns", "PF,C",  "New York Knicks"),
    @("Trey Murphy III",    "SF,PF", "New Orleans Pelicans"),
    @("Daniel Gafford",     "PF,C",  "Dallas Mavericks"),
    @("Mark Williams",      "C",     "Charlotte Hornets"),
    @("Jalen Duren",        "C",     "Detroit Pistons"),
    @("Keegan Murray",      "SF,PF", "Sacramento Kings"),
    @("Tyrese Maxey",       "PG,SG", "Philadelphia 76ers"),
    @("Tyrese Haliburton",  "PG,SG", "Indiana Pacers"),
    @("Jarrett Allen",      "C",     "Cleveland Cavaliers"),
    @("Jalen Johnson",      "SF,PF", "Atlanta Hawks"),
    @("Franz Wagner",       "SF,PF", "Orlando Magic"),
    @("Malcolm Brogdon",    "PG,SG", "Washington Wizards")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
}
